$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, copying the header style (bold/border/centered)
# from the existing "sum" header cell (G1) so it reuses the same style record.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill H2:H12 with the "Save" indicator values (1 for a save, 0 otherwise).
$saveValues = @(0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
